$p = $ppt.ActivePresentation

# --- Slide 7: "Informe sobre los partes de asignatura de docentes" ---
# Text box "CuadroTexto 6" is shape index 3.
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(3)
$tr7 = $shp7.TextFrame.TextRange

# "Observaciones destacables: {{OBSERVACIONES_DESTACABLES}}" ->
# "Observaciones destacables: {{OBSERVACIONES_DESTACABLES_DOC}}"
$full7 = $tr7.Text
$oldObs = "Observaciones destacables: {{OBSERVACIONES_DESTACABLES}}"
$newObs = "Observaciones destacables: {{OBSERVACIONES_DESTACABLES_DOC}}"
$idxObs = $full7.IndexOf($oldObs)
if ($idxObs -ge 0) {
    $rngObs = $tr7.Characters($idxObs + 1, $oldObs.Length)
    $rngObs.Text = $newObs
}

# "Comentarios destacables: {{COMENTARIOS_DESTACABLES}}" ->
# "Comentarios destacables: {{COMENTARIOS" + "_DESTACABLES_DOC}}"
$full7b = $tr7.Text
$splitAnchor = "COMENTARIOS_DESTACABLES}}"
$idxCom = $full7b.LastIndexOf($splitAnchor)
if ($idxCom -ge 0) {
    $prefixLen = "COMENTARIOS".Length
    $tailStart = $idxCom + $prefixLen
    $tailLen = $splitAnchor.Length - $prefixLen
    $rngTail = $tr7.Characters($tailStart + 1, $tailLen)
    $rngTail.Text = "_DESTACABLES_DOC}}"
}

# --- Slide 8: "Informe sobre los partes de asignatura de delegados" ---
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(3)
$tr8 = $shp8.TextFrame.TextRange

# "Observaciones destacables: {{OBSERVACIONES_DESTACABLES}}" ->
# "Observaciones destacables: {{OBSERVACIONES_DESTACABLES_DEL}}"
$full8 = $tr8.Text
$oldObs8 = "Observaciones destacables: {{OBSERVACIONES_DESTACABLES}}"
$newObs8 = "Observaciones destacables: {{OBSERVACIONES_DESTACABLES_DEL}}"
$idxObs8 = $full8.IndexOf($oldObs8)
if ($idxObs8 -ge 0) {
    $rngObs8 = $tr8.Characters($idxObs8 + 1, $oldObs8.Length)
    $rngObs8.Text = $newObs8
}

# "Comentarios destacables: {{COMENTARIOS_DESTACABLES}}" ->
# "Comentarios destacables: {{COMENTARIOS_DESTACABLES_DEL}}"
$full8b = $tr8.Text
$oldCom8 = "Comentarios destacables: {{COMENTARIOS_DESTACABLES}}"
$newCom8 = "Comentarios destacables: {{COMENTARIOS_DESTACABLES_DEL}}"
$idxCom8 = $full8b.IndexOf($oldCom8)
if ($idxCom8 -ge 0) {
    $rngCom8 = $tr8.Characters($idxCom8 + 1, $oldCom8.Length)
    $rngCom8.Text = $newCom8
}
